$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.29779169861954
$ws.Range("C2").Value = 8.272530423114503
$ws.Range("D2").Value = 13.52983058904382
$ws.Range("E2").Value = 13.90758626613974
$ws.Range("G2").Value = 3.722759816181884
$ws.Range("J2").Value = 8.525373922252752
$ws.Range("L2").Value = 12.25760396085246
$ws.Range("M2").Value = 18.77074067939193
$ws.Range("O2").Value = 33.78363315571368

$ws.Range("B3").Value = 19.96064240996766
$ws.Range("C3").Value = 8.023605634309186
$ws.Range("D3").Value = 13.55077148724115
$ws.Range("E3").Value = 13.94595018123578
$ws.Range("G3").Value = 3.725484214084565
$ws.Range("J3").Value = 8.525381100959216
$ws.Range("L3").Value = 12.26438464911869
$ws.Range("M3").Value = 18.7028360213595
$ws.Range("O3").Value = 33.88382460644006

$ws.Range("B4").Value = 19.75429392280465
$ws.Range("C4").Value = 7.865680470976323
$ws.Range("D4").Value = 13.56590951430752
$ws.Range("E4").Value = 13.9708779051951
$ws.Range("G4").Value = 3.727245587062147
$ws.Range("J4").Value = 8.525487238656009
$ws.Range("L4").Value = 12.27003308100755
$ws.Range("M4").Value = 18.66374364563452
$ws.Range("O4").Value = 33.9528645194572

$ws.Range("B5").Value = 19.67047774231453
$ws.Range("C5").Value = 7.800110681707038
$ws.Range("D5").Value = 13.57265150990957
$ws.Range("E5").Value = 13.98138196061439
$ws.Range("G5").Value = 3.727985710457792
$ws.Range("J5").Value = 8.525556177874128
$ws.Range("L5").Value = 12.27270899521427
$ws.Range("M5").Value = 18.64847803203463
$ws.Range("O5").Value = 33.98288495532304

$ws.Range("B6").Value = 19.65657969790023
$ws.Range("C6").Value = 7.789151543044562
$ws.Range("D6").Value = 13.5738056192035
$ws.Range("E6").Value = 13.98314705962624
$ws.Range("G6").Value = 3.728109959451335
$ws.Range("J6").Value = 8.525569179756426
$ws.Range("L6").Value = 12.27317594460155
$ws.Range("M6").Value = 18.64598364402306
$ws.Range("O6").Value = 33.987983581708

$ws.Range("B7").Value = 19.75316230464719
$ws.Range("C7").Value = 7.86480099986712
$ws.Range("D7").Value = 13.56599811905775
$ws.Range("E7").Value = 13.97101816537823
$ws.Range("G7").Value = 3.727255478032916
$ws.Range("J7").Value = 8.525488064247051
$ws.Range("L7").Value = 12.27006765359855
$ws.Range("M7").Value = 18.66353506263376
$ws.Range("O7").Value = 33.95326175605776

$ws.Range("B8").Value = 20.18146834464057
$ws.Range("C8").Value = 8.187793327128775
$ws.Range("D8").Value = 13.53657762159376
$ws.Range("E8").Value = 13.92052982560722
$ws.Range("G8").Value = 3.723680848696202
$ws.Range("J8").Value = 8.52535534660432
$ws.Range("L8").Value = 12.25963415575909
$ws.Range("M8").Value = 18.7467938559028
$ws.Range("O8").Value = 33.8166148349852

$ws.Range("B9").Value = 21.02204417037282
$ws.Range("C9").Value = 8.778295192770647
$ws.Range("D9").Value = 13.49698545951971
$ws.Range("E9").Value = 13.8323742942494
$ws.Range("G9").Value = 3.71737044520378
$ws.Range("J9").Value = 8.525897829521051
$ws.Range("L9").Value = 12.25092513073666
$ws.Range("M9").Value = 18.93020791312015
$ws.Range("O9").Value = 33.60856504825626

$ws.Range("B10").Value = 21.63402602293745
$ws.Range("C10").Value = 9.182978177324044
$ws.Range("D10").Value = 13.47894001266389
$ws.Range("E10").Value = 13.77417281255399
$ws.Range("G10").Value = 3.713155797789832
$ws.Range("J10").Value = 8.526780143971969
$ws.Range("L10").Value = 12.25164483406734
$ws.Range("M10").Value = 19.07654076990814
$ws.Range("O10").Value = 33.49254013057749

$ws.Range("B11").Value = 21.91000681680409
$ws.Range("C11").Value = 9.360214479861479
$ws.Range("D11").Value = 13.47312810590642
$ws.Range("E11").Value = 13.74911087742174
$ws.Range("G11").Value = 3.711328973753039
$ws.Range("J11").Value = 8.527285410178383
$ws.Range("L11").Value = 12.25350682077755
$ws.Range("M11").Value = 19.14546460581499
$ws.Range("O11").Value = 33.44781448835115

$ws.Range("B12").Value = 22.01407056344775
$ws.Range("C12").Value = 9.426305227210747
$ws.Range("D12").Value = 13.47127171888288
$ws.Range("E12").Value = 13.73982316360549
$ws.Range("G12").Value = 3.710650130556998
$ws.Range("J12").Value = 8.527491579535708
$ws.Range("L12").Value = 12.2544314867848
$ws.Range("M12").Value = 19.17188866683951
$ws.Range("O12").Value = 33.43204076932448

$ws.Range("B13").Value = 21.99167983019979
$ws.Range("C13").Value = 9.412117552834111
$ws.Range("D13").Value = 13.47165621036834
$ws.Range("E13").Value = 13.74181443443532
$ws.Range("G13").Value = 3.710795757385241
$ws.Range("J13").Value = 8.527446519094866
$ws.Range("L13").Value = 12.25422259624388
$ws.Range("M13").Value = 19.16618359732442
$ws.Range("O13").Value = 33.43538613030118

$ws.Range("B14").Value = 21.91857764732746
$ws.Range("C14").Value = 9.365672550953729
$ws.Range("D14").Value = 13.47296847821677
$ws.Range("E14").Value = 13.74834271304016
$ws.Range("G14").Value = 3.711272866067564
$ws.Range("J14").Value = 8.527302075115269
$ws.Range("L14").Value = 12.25357849963633
$ws.Range("M14").Value = 19.14763211972308
$ws.Range("O14").Value = 33.44649344295499

$ws.Range("B15").Value = 21.87373979234048
$ws.Range("C15").Value = 9.337089108809792
$ws.Range("D15").Value = 13.47381712934543
$ws.Range("E15").Value = 13.75236784975799
$ws.Range("G15").Value = 3.711566791422136
$ws.Range("J15").Value = 8.52721552800781
$ws.Range("L15").Value = 12.25321253263752
$ws.Range("M15").Value = 19.1363105666342
$ws.Range("O15").Value = 33.45344857094373

$ws.Range("B16").Value = 21.61593306598499
$ws.Range("C16").Value = 9.171253979486224
$ws.Range("D16").Value = 13.47936805180296
$ws.Range("E16").Value = 13.77583906946745
$ws.Range("G16").Value = 3.713276998949999
$ws.Range("J16").Value = 8.526749205387169
$ws.Range("L16").Value = 12.25155393648437
$ws.Range("M16").Value = 19.07208260076078
$ws.Range("O16").Value = 33.49562554242954

$ws.Range("B17").Value = 21.45708846009343
$ws.Range("C17").Value = 9.067734394673998
$ws.Range("D17").Value = 13.48338717615963
$ws.Range("E17").Value = 13.79059964005835
$ws.Range("G17").Value = 3.714349269988806
$ws.Range("J17").Value = 8.526489667737978
$ws.Range("L17").Value = 12.25092872899168
$ws.Range("M17").Value = 19.03327399030605
$ws.Range("O17").Value = 33.5235661563048

$ws.Range("B18").Value = 21.36550368985703
$ws.Range("C18").Value = 9.00754937664734
$ws.Range("D18").Value = 13.48592451450659
$ws.Range("E18").Value = 13.79922268123709
$ws.Range("G18").Value = 3.714974528729838
$ws.Range("J18").Value = 8.526350174018212
$ws.Range("L18").Value = 12.25071363983433
$ws.Range("M18").Value = 19.01117500053371
$ws.Range("O18").Value = 33.54039474287765

$ws.Range("B19").Value = 21.33445964176148
$ws.Range("C19").Value = 8.987062536922052
$ws.Range("D19").Value = 13.48682237355992
$ws.Range("E19").Value = 13.80216518409853
$ws.Range("G19").Value = 3.715187695498894
$ws.Range("J19").Value = 8.526304627690461
$ws.Range("L19").Value = 12.25066566693774
$ws.Range("M19").Value = 19.0037313520185
$ws.Range("O19").Value = 33.54622263360066

$ws.Range("B20").Value = 21.47402138775228
$ws.Range("C20").Value = 9.078821105373198
$ws.Range("D20").Value = 13.48293598150998
$ws.Range("E20").Value = 13.78901457629322
$ws.Range("G20").Value = 3.714234243955438
$ws.Range("J20").Value = 8.526516283960666
$ws.Range("L20").Value = 12.25098033449422
$ws.Range("M20").Value = 19.03738229014618
$ws.Range("O20").Value = 33.5205133572665

$ws.Range("B21").Value = 21.9400623386301
$ws.Range("C21").Value = 9.379342674605324
$ws.Range("D21").Value = 13.47257368740872
$ws.Range("E21").Value = 13.74641970347741
$ws.Range("G21").Value = 3.711132377172153
$ws.Range("J21").Value = 8.527344099996547
$ws.Range("L21").Value = 12.25376173645772
$ws.Range("M21").Value = 19.15307246200033
$ws.Range("O21").Value = 33.44319935933185

$ws.Range("B22").Value = 22.24201603953141
$ws.Range("C22").Value = 9.56976261520332
$ws.Range("D22").Value = 13.46780888475229
$ws.Range("E22").Value = 13.71976266354737
$ws.Range("G22").Value = 3.709180493397064
$ws.Range("J22").Value = 8.52797155322455
$ws.Range("L22").Value = 12.25685881150987
$ws.Range("M22").Value = 19.23056464502161
$ws.Range("O22").Value = 33.39945046548645

$ws.Range("B23").Value = 22.08112936951427
$ws.Range("C23").Value = 9.468691562590843
$ws.Range("D23").Value = 13.47016836847009
$ws.Range("E23").Value = 13.73388216743312
$ws.Range("G23").Value = 3.710215377487863
$ws.Range("J23").Value = 8.527628795814357
$ws.Range("L23").Value = 12.25508917491523
$ws.Range("M23").Value = 19.18903838731813
$ws.Range("O23").Value = 33.42217822041741

$ws.Range("B24").Value = 21.46636682252891
$ws.Range("C24").Value = 9.073810888146207
$ws.Range("D24").Value = 13.48313926048045
$ws.Range("E24").Value = 13.78973075684482
$ws.Range("G24").Value = 3.71428621982061
$ws.Range("J24").Value = 8.526504220487102
$ws.Range("L24").Value = 12.25095655392245
$ws.Range("M24").Value = 19.03552426443377
$ws.Range("O24").Value = 33.52189114455551

$ws.Range("B25").Value = 20.79520431593489
$ws.Range("C25").Value = 8.623471887657201
$ws.Range("D25").Value = 13.50575706709566
$ws.Range("E25").Value = 13.8550660491403
$ws.Range("G25").Value = 3.719003193993786
$ws.Range("J25").Value = 8.525665717220859
$ws.Range("L25").Value = 12.25202752840804
$ws.Range("M25").Value = 18.87850313569844
$ws.Range("O25").Value = 33.65840138682113
